$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$data = @{
    3  = @(1279, 2319, 5139, 10500, 14800, 14500)
    8  = @(15600, 16000, 28200, 72600, 134000, 236000)
    13 = @(10800, 16400, 19700, 21300, 21800, 22200)
    18 = @(358000, 428000, 525000, 589000, 618000, 680000)
    23 = @(22300, 41800, 42000, 41400, 41700, 42000)
    28 = @(110000, 1008000, 1074000, 1216000, 1226000, 1224000)
    33 = @(10200, 12300, 13300, 13500, 14100, 14600)
    38 = @(302000, 363000, 406000, 427000, 447000, 459000)
}

$cols = @("B", "C", "D", "E", "F", "G")

foreach ($row in $data.Keys) {
    $values = $data[$row]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $cellRef = "$($cols[$i])$row"
        $ws.Range($cellRef).Value = $values[$i]
    }
}
